$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 204-213 (values shift as the weekly price sheet refreshes) ---

# Row 204
$ws.Cells.Item(204, "D").Value = 44516
$ws.Cells.Item(204, "K").Value = "Early Burlat"
$ws.Cells.Item(204, "L").Value = "Primera"
$ws.Cells.Item(204, "M").Value = 35
$ws.Cells.Item(204, "N").Value = 12500
$ws.Cells.Item(204, "O").Value = 12500
$ws.Cells.Item(204, "P").Value = 12500
$ws.Cells.Item(204, "Q").Value = "`$/bandeja 5 kilos"
$ws.Cells.Item(204, "R").Value = "Provincia de Curicó"
$ws.Cells.Item(204, "S").Value = 2500
$ws.Cells.Item(204, "T").Value = 5

# Row 205
$ws.Cells.Item(205, "D").Value = 44516
$ws.Cells.Item(205, "K").Value = "Early Burlat"
$ws.Cells.Item(205, "L").Value = "Segunda"
$ws.Cells.Item(205, "M").Value = 30
$ws.Cells.Item(205, "N").Value = 10000
$ws.Cells.Item(205, "O").Value = 10000
$ws.Cells.Item(205, "P").Value = 10000
$ws.Cells.Item(205, "Q").Value = "`$/bandeja 5 kilos"
$ws.Cells.Item(205, "R").Value = "Provincia de Curicó"
$ws.Cells.Item(205, "S").Value = 2000
$ws.Cells.Item(205, "T").Value = 5

# Row 206
$ws.Cells.Item(206, "D").Value = 44217
$ws.Cells.Item(206, "K").Value = "Santina"
$ws.Cells.Item(206, "L").Value = "Especial"
$ws.Cells.Item(206, "M").Value = 56
$ws.Cells.Item(206, "N").Value = 12000
$ws.Cells.Item(206, "O").Value = 12000
$ws.Cells.Item(206, "P").Value = 12000
$ws.Cells.Item(206, "Q").Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(206, "R").Value = "Provincia de Curicó"
$ws.Cells.Item(206, "S").Value = 1200
$ws.Cells.Item(206, "T").Value = 10

# Row 207
$ws.Cells.Item(207, "D").Value = 44217
$ws.Cells.Item(207, "K").Value = "Santina"
$ws.Cells.Item(207, "L").Value = "Primera"
$ws.Cells.Item(207, "M").Value = 60
$ws.Cells.Item(207, "N").Value = 10000
$ws.Cells.Item(207, "O").Value = 10000
$ws.Cells.Item(207, "P").Value = 10000
$ws.Cells.Item(207, "Q").Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(207, "R").Value = "Provincia de Curicó"
$ws.Cells.Item(207, "S").Value = 1000
$ws.Cells.Item(207, "T").Value = 10

# Row 208
$ws.Cells.Item(208, "D").Value = 44217
$ws.Cells.Item(208, "K").Value = "Santina"
$ws.Cells.Item(208, "L").Value = "Segunda"
$ws.Cells.Item(208, "M").Value = 50
$ws.Cells.Item(208, "N").Value = 8000
$ws.Cells.Item(208, "O").Value = 8000
$ws.Cells.Item(208, "P").Value = 8000
$ws.Cells.Item(208, "Q").Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(208, "R").Value = "Provincia de Curicó"
$ws.Cells.Item(208, "S").Value = 800
$ws.Cells.Item(208, "T").Value = 10

# Row 209
$ws.Cells.Item(209, "D").Value = 44509
$ws.Cells.Item(209, "K").Value = "Early Burlat"
$ws.Cells.Item(209, "L").Value = "Segunda"
$ws.Cells.Item(209, "M").Value = 36
$ws.Cells.Item(209, "N").Value = 40000
$ws.Cells.Item(209, "O").Value = 40000
$ws.Cells.Item(209, "P").Value = 40000
$ws.Cells.Item(209, "Q").Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(209, "R").Value = "Provincia de Curicó"
$ws.Cells.Item(209, "S").Value = 4000
$ws.Cells.Item(209, "T").Value = 10

# Row 210
$ws.Cells.Item(210, "D").Value = 44179
$ws.Cells.Item(210, "K").Value = "Lapins"
$ws.Cells.Item(210, "L").Value = "Especial"
$ws.Cells.Item(210, "M").Value = 75
$ws.Cells.Item(210, "N").Value = 12000
$ws.Cells.Item(210, "O").Value = 12000
$ws.Cells.Item(210, "P").Value = 12000
$ws.Cells.Item(210, "Q").Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(210, "R").Value = "Provincia de Curicó"
$ws.Cells.Item(210, "S").Value = 1200
$ws.Cells.Item(210, "T").Value = 10

# Row 211
$ws.Cells.Item(211, "D").Value = 44179
$ws.Cells.Item(211, "K").Value = "Lapins"
$ws.Cells.Item(211, "L").Value = "Primera"
$ws.Cells.Item(211, "M").Value = 78
$ws.Cells.Item(211, "N").Value = 10000
$ws.Cells.Item(211, "O").Value = 10000
$ws.Cells.Item(211, "P").Value = 10000
$ws.Cells.Item(211, "Q").Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(211, "R").Value = "Provincia de Curicó"
$ws.Cells.Item(211, "S").Value = 1000
$ws.Cells.Item(211, "T").Value = 10

# Row 212
$ws.Cells.Item(212, "D").Value = 44179
$ws.Cells.Item(212, "K").Value = "Lapins"
$ws.Cells.Item(212, "L").Value = "Segunda"
$ws.Cells.Item(212, "M").Value = 70
$ws.Cells.Item(212, "N").Value = 9000
$ws.Cells.Item(212, "O").Value = 9000
$ws.Cells.Item(212, "P").Value = 9000
$ws.Cells.Item(212, "Q").Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(212, "R").Value = "Provincia de Curicó"
$ws.Cells.Item(212, "S").Value = 900
$ws.Cells.Item(212, "T").Value = 10

# Row 213
$ws.Cells.Item(213, "D").Value = 44179
$ws.Cells.Item(213, "K").Value = "Rainier"
$ws.Cells.Item(213, "L").Value = "Especial"
$ws.Cells.Item(213, "M").Value = 68
$ws.Cells.Item(213, "N").Value = 16000
$ws.Cells.Item(213, "O").Value = 16000
$ws.Cells.Item(213, "P").Value = 16000
$ws.Cells.Item(213, "Q").Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(213, "R").Value = "Región del Maule"
$ws.Cells.Item(213, "S").Value = 1600
$ws.Cells.Item(213, "T").Value = 10

# Row 214
# New row 214 -- fill fixed columns first, then set the date columns number format
# to match the rest of the column (so it renders/serialises as a date, not a raw serial).
$ws.Cells.Item(214, "A").Value = 3
$ws.Cells.Item(214, "B").Value = "Femacal de La Calera"
$ws.Cells.Item(214, "C").Value = "Coquimbo"
$ws.Cells.Item(214, "E").Value = 5
$ws.Cells.Item(214, "F").Value = "Fruta"
$ws.Cells.Item(214, "G").Value = 100103
$ws.Cells.Item(214, "H").Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(214, "I").Value = 100103001
$ws.Cells.Item(214, "J").Value = "Cereza"
$ws.Cells.Item(214, "D").Value = 44179
$ws.Cells.Item(214, "D").NumberFormat = $ws.Cells.Item(213, "D").NumberFormat
$ws.Cells.Item(214, "K").Value = "Rainier"
$ws.Cells.Item(214, "L").Value = "Primera"
$ws.Cells.Item(214, "M").Value = 78
$ws.Cells.Item(214, "N").Value = 14000
$ws.Cells.Item(214, "O").Value = 14000
$ws.Cells.Item(214, "P").Value = 14000
$ws.Cells.Item(214, "Q").Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(214, "R").Value = "Región del Maule"
$ws.Cells.Item(214, "S").Value = 1400
$ws.Cells.Item(214, "T").Value = 10

# Row 215
# New row 215 -- fill fixed columns first, then set the date columns number format
# to match the rest of the column (so it renders/serialises as a date, not a raw serial).
$ws.Cells.Item(215, "A").Value = 3
$ws.Cells.Item(215, "B").Value = "Femacal de La Calera"
$ws.Cells.Item(215, "C").Value = "Coquimbo"
$ws.Cells.Item(215, "E").Value = 5
$ws.Cells.Item(215, "F").Value = "Fruta"
$ws.Cells.Item(215, "G").Value = 100103
$ws.Cells.Item(215, "H").Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(215, "I").Value = 100103001
$ws.Cells.Item(215, "J").Value = "Cereza"
$ws.Cells.Item(215, "D").Value = 44179
$ws.Cells.Item(215, "D").NumberFormat = $ws.Cells.Item(213, "D").NumberFormat
$ws.Cells.Item(215, "K").Value = "Rainier"
$ws.Cells.Item(215, "L").Value = "Segunda"
$ws.Cells.Item(215, "M").Value = 70
$ws.Cells.Item(215, "N").Value = 13000
$ws.Cells.Item(215, "O").Value = 13000
$ws.Cells.Item(215, "P").Value = 13000
$ws.Cells.Item(215, "Q").Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(215, "R").Value = "Región del Maule"
$ws.Cells.Item(215, "S").Value = 1300
$ws.Cells.Item(215, "T").Value = 10
